# Auto-generated Excel COM-interop script applying the diff to Gilgamesh_Profits workbook
# Updates cached values in columns H-N across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 829.93335
$ws.Range("J29").Value = 882.0714
$ws.Range("L29").Value = 2646.2142
$ws.Range("N29").Value = -3208.2142

$ws.Range("H61").Value = 273
$ws.Range("I61").Value = 273
$ws.Range("K61").Value = 819
$ws.Range("M61").Value = -647

$ws.Range("H92").Value = 1153
$ws.Range("I92").Value = 1153
$ws.Range("K92").Value = 1153
$ws.Range("M92").Value = 95

$ws.Range("H107").Value = 603.6111
$ws.Range("I107").Value = 521.4706
$ws.Range("K107").Value = 521.4706
$ws.Range("M107").Value = 1398.5294

$ws.Range("H127").Value = 984.1667
$ws.Range("I127").Value = 984.1667
$ws.Range("K127").Value = 2952.5001
$ws.Range("M127").Value = 2007.4999

$ws.Range("H141").Value = 3535.2
$ws.Range("I141").Value = 3377.4167
$ws.Range("J141").Value = 4166.3335
$ws.Range("K141").Value = 10132.2501
$ws.Range("L141").Value = 12499.0005
$ws.Range("M141").Value = -4952.250100000001
$ws.Range("N141").Value = -22859.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1568969.1
$ws.Range("I32").Value = 742046.1
$ws.Range("K32").Value = 742046.1
$ws.Range("M32").Value = -741759.1

$ws.Range("H74").Value = 1205.225
$ws.Range("I74").Value = 727.96875
$ws.Range("K74").Value = 727.96875
$ws.Range("M74").Value = 146.03125

$ws.Range("H77").Value = 1205.225
$ws.Range("I77").Value = 727.96875
$ws.Range("K77").Value = 3639.84375
$ws.Range("M77").Value = 728.15625

$ws.Range("H102").Value = 2192.4211
$ws.Range("I102").Value = 1610.4667
$ws.Range("J102").Value = 4374.75
$ws.Range("K102").Value = 1610.4667
$ws.Range("L102").Value = 4374.75
$ws.Range("M102").Value = 11.53330000000005
$ws.Range("N102").Value = -7618.75

$ws.Range("H122").Value = 2064.1904
$ws.Range("J122").Value = 1023.8571
$ws.Range("L122").Value = 3071.5713
$ws.Range("N122").Value = -7971.5713

$ws.Range("H132").Value = 2154.5
$ws.Range("I132").Value = 1667.6111
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 5002.8333
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -2472.8333
$ws.Range("N132").Value = -14810

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 166688670
$ws.Range("I20").Value = 250030000
$ws.Range("J20").Value = 5999
$ws.Range("K20").Value = 250030000
$ws.Range("L20").Value = 5999
$ws.Range("M20").Value = -250029753
$ws.Range("N20").Value = -6493

$ws.Range("H64").Value = 1530.1
$ws.Range("J64").Value = 2247.5
$ws.Range("L64").Value = 2247.5
$ws.Range("N64").Value = -2697.5

$ws.Range("H67").Value = 1530.1
$ws.Range("J67").Value = 2247.5
$ws.Range("L67").Value = 2247.5
$ws.Range("N67").Value = -3807.5

$ws.Range("H97").Value = 29999.75
$ws.Range("J97").Value = 29999.75
$ws.Range("L97").Value = 29999.75
$ws.Range("N97").Value = -31981.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 115.15
$ws.Range("I7").Value = 68.84614999999999
$ws.Range("J7").Value = 201.14285
$ws.Range("K7").Value = 68.84614999999999
$ws.Range("L7").Value = 201.14285
$ws.Range("M7").Value = 44.15385000000001
$ws.Range("N7").Value = -427.14285

$ws.Range("H132").Value = 4199.4517
$ws.Range("I132").Value = 3532.5715
$ws.Range("K132").Value = 10597.7145
$ws.Range("M132").Value = -8067.7145

$ws.Range("H134").Value = 3752.2163
$ws.Range("I134").Value = 3769.6897
$ws.Range("K134").Value = 11309.0691
$ws.Range("M134").Value = -8774.069100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 5375
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 5375
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 16125
$ws.Range("N54").Value = -17243
$ws.Range("M54").Value = $null

$ws.Range("H121").Value = 5055408.5
$ws.Range("I121").Value = 11111518
$ws.Range("J121").Value = 100409.63
$ws.Range("K121").Value = 33334554
$ws.Range("L121").Value = 301228.89
$ws.Range("M121").Value = -33333244
$ws.Range("N121").Value = -303848.89

$ws.Range("H124").Value = 8000
$ws.Range("I124").Value = 8000
$ws.Range("K124").Value = 24000
$ws.Range("M124").Value = -19090

$ws.Range("H131").Value = 3107919.5
$ws.Range("J131").Value = 4904036.5
$ws.Range("L131").Value = 14712109.5
$ws.Range("N131").Value = -14722189.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 4952.5
$ws.Range("J6").Value = 4952.5
$ws.Range("L6").Value = 4952.5
$ws.Range("N6").Value = -5178.5

$ws.Range("H16").Value = 4952.5
$ws.Range("J16").Value = 4952.5
$ws.Range("L16").Value = 4952.5
$ws.Range("N16").Value = -5452.5

$ws.Range("H52").Value = 26500
$ws.Range("J52").Value = 28000
$ws.Range("L52").Value = 28000
$ws.Range("N52").Value = -28518

$ws.Range("H70").Value = 5317.636
$ws.Range("I70").Value = 4277.1113
$ws.Range("K70").Value = 4277.1113
$ws.Range("M70").Value = -4007.1113

$ws.Range("H73").Value = 5317.636
$ws.Range("I73").Value = 4277.1113
$ws.Range("K73").Value = 4277.1113
$ws.Range("M73").Value = -3341.1113

$ws.Range("H107").Value = 1469.5883
$ws.Range("I107").Value = 370.875
$ws.Range("K107").Value = 370.875
$ws.Range("M107").Value = 1549.125

$ws.Range("H132").Value = 2372.5518
$ws.Range("I132").Value = 2300.25
$ws.Range("J132").Value = 2533.2222
$ws.Range("K132").Value = 6900.75
$ws.Range("L132").Value = 7599.6666
$ws.Range("M132").Value = -4370.75
$ws.Range("N132").Value = -12659.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1734.4
$ws.Range("I22").Value = 1168
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 1168
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -873
$ws.Range("N22").Value = -4590

$ws.Range("H27").Value = 1734.4
$ws.Range("I27").Value = 1168
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 1168
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = -1061
$ws.Range("N27").Value = -4214

$ws.Range("H40").Value = 23058.354
$ws.Range("I40").Value = 24187
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 24187
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -24051
$ws.Range("N40").Value = -5272

$ws.Range("H119").Value = 99000
$ws.Range("J119").Value = 99000
$ws.Range("L119").Value = 99000
$ws.Range("N119").Value = -108676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13160015
$ws.Range("I122").Value = 2049.1875
$ws.Range("J122").Value = 83335830
$ws.Range("K122").Value = 6147.5625
$ws.Range("L122").Value = 250007490
$ws.Range("M122").Value = -3697.5625
$ws.Range("N122").Value = -250012390

$ws.Range("H132").Value = 2835.1365
$ws.Range("I132").Value = 2835.1365
$ws.Range("K132").Value = 8505.4095
$ws.Range("M132").Value = -5975.4095
